# This workbook's rows 3-7 (data rows for "Fruta, Feria Lagunitas de Puerto
# Montt - Tuna") get reordered/updated: columns D, L, M, N, O, P, R, S change
# per-row so that each row ends up holding a different combination of values
# than it started with (effectively a re-sort of the weekly records).
#
# Target values per row (columns: D=Fecha, L=Calidad, M=Volumen,
# N=Precio minimo, O=Precio maximo, P=Precio promedio ponderado,
# R=Origen, S=Precio $/Kg):

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    3 = @{ D = 44253; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana"; S = 806 }
    4 = @{ D = 44250; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana"; S = 806 }
    5 = @{ D = 44257; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana"; S = 806 }
    6 = @{ D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco"; S = 889 }
    7 = @{ D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí"; S = 944 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
